{"js": "// Gibson - Twilio Executive Summary: apply the author's edits.\n//\n// 1) \"...roughly 1% mortality rate). For cases...\" ->\n//    \"...roughly 1% mortality rate) in Hamilton County, Ohio (containing\n//    Cincinnati). For cases...\"\n// 2) \"...a negative effect on mortality rates. The data shows a moderate \"\n//    is retyped (no wording change, runs collapse to one).\n// 3) \"Generalize model and expand analysis...\" is retyped (no wording\n//    change, runs collapse to one).\n\nasync function replaceWithSelf(body, text) {\n  const results = body.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    return false;\n  }\n  results.items[0].insertText(text, Word.InsertLocation.replace);\n  await context.sync();\n  return true;\n}\n\nconst body = context.document.body;\n\n// 1) Insert the new geography qualifier right after the mortality-rate\n// parenthetical and before the trailing hospitalization sentence.\nconst anchor = \"140k cases with ~1700 deaths (roughly 1% mortality rate)\";\nconst anchorResults = body.search(anchor, { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Could not find the mortality-rate sentence to edit.\");\n}\n\nanchorResults.items[0].insertText(\n  \" in Hamilton County, Ohio (containing Cincinnati)\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// 2) No wording change, just re-typed as a single run in the source edit.\nawait replaceWithSelf(\n  body,\n  \"I investigated to find the correlation between vaccination and mortality rates to confirm that vaccinations have a negative effect on mortality rates. The data shows a moderate \"\n);\n\n// 3) No wording change, just re-typed as a single run in the source edit.\nawait replaceWithSelf(\n  body,\n  \"Generalize model and expand analysis to include additional geographies to confirm findings\"\n);\n", "ps1": "# Gibson - Twilio Executive Summary: apply the author's edits.\n#\n# 1) \"...roughly 1% mortality rate). For cases...\" ->\n#    \"...roughly 1% mortality rate) in Hamilton County, Ohio (containing\n#    Cincinnati). For cases...\"\n# 2) \"...a negative effect on mortality rates. The data shows a moderate \"\n#    is retyped (no wording change, runs collapse to one).\n# 3) \"Generalize model and expand analysis...\" is retyped (no wording\n#    change, runs collapse to one).\n\n$d = $word.ActiveDocument\n\n# 1) Insert the new geography qualifier right after the mortality-rate\n# parenthetical and before the trailing hospitalization sentence.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"140k cases with ~1700 deaths (roughly 1% mortality rate)\"\n$find.MatchCase = $true\n$found = $find.Execute()\nif ($found) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\" in Hamilton County, Ohio (containing Cincinnati)\")\n}\n\n# 2) No wording change, just re-typed as a single run in the source edit.\n$text2 = \"I investigated to find the correlation between vaccination and mortality rates to confirm that vaccinations have a negative effect on mortality rates. The data shows a moderate \"\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = $text2\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2.Text = $text2\n}\n\n# 3) No wording change, just re-typed as a single run in the source edit.\n$text3 = \"Generalize model and expand analysis to include additional geographies to confirm findings\"\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$find3.Text = $text3\n$find3.MatchCase = $true\n$found3 = $find3.Execute()\nif ($found3) {\n    $rng3.Text = $text3\n}\n"}
